# Add two new rows (17 and 18) of AV data below the existing table,
# matching the "create form for AV PEG 600DO-IK" commit.
#
# Columns: A=LOT, B=Step, C=Waktu, D=Operator QC, E=Reaksi (C),
#          F=Berat Sample (gr), G=Jumlah Titran (mL), H=Faktor Buret,
#          I=Faktor NaOH, J=AV, K=Instruksi
#
# A, C, D, E, K must be stored as text even though several of the values
# look numeric (e.g. "12", "123", "4123"), matching the source which used
# inlineStr for those columns. B, F, G, H, I, J stay as real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{
        Row = 17
        A = "4123"
        B = 1
        C = "14:32"
        D = "12"
        E = "123"
        F = 12
        G = 1
        H = 12
        I = 12
        J = 807.84
        K = "Hubungi atasan"
    },
    @{
        Row = 18
        A = "12"
        B = 1
        C = "14:35"
        D = "12"
        E = "123"
        F = 12
        G = 12
        H = 12
        I = 12
        J = 9694.08
        K = "Hubungi atasan"
    }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Force the text-typed columns to be stored as text before assigning
    # their values, so number-looking strings like "12" are not coerced
    # into numeric cells.
    foreach ($col in @("A", "C", "D", "E", "K")) {
        $ws.Range("$col$rowNum").NumberFormat = "@"
    }

    $ws.Range("A$rowNum").Value = $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("F$rowNum").Value = $r.F
    $ws.Range("G$rowNum").Value = $r.G
    $ws.Range("H$rowNum").Value = $r.H
    $ws.Range("I$rowNum").Value = $r.I
    $ws.Range("J$rowNum").Value = $r.J
    $ws.Range("K$rowNum").Value = $r.K

    # Clear the style index we picked up from the NumberFormat change so
    # the new cells stay unstyled, just like the rest of the data rows.
    $ws.Range("A$($rowNum):K$($rowNum)").Style = "Normal"
}

Write-Host "Added rows 17-18; used range is now" $ws.UsedRange.Address()
